$d = $word.ActiveDocument

# Grab the whole document content as a WordOpenXML package so we can do
# precise, surgical text-level surgery on the underlying OOXML and then
# write the whole story back in one shot via InsertXML.
$rng = $d.Content
$xml = $rng.WordOpenXML

# --- Change 1 --------------------------------------------------------
# "a easily way" -> split the single run into three runs with
# proofErr gramStart/gramEnd bracketing the lone "a", matching Word's
# grammar-checker markup for the "a easily" flag.
$old1 = '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>In this solution lets just say the man can wait and get on a bigger boat and in this case he is able to carry the bag of seeds, parrot and cat all the way across the river to the other side. This will give him a easily way to transport all three items at once.</w:t></w:r>'
$new1 = '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">In this solution lets just say the man can wait and get on a bigger boat and in this case he is able to carry the bag of seeds, parrot and cat all the way across the river to the other side. This will give him </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> easily way to transport all three items at once.</w:t></w:r>'

if (-not $xml.Contains($old1)) {
    throw "edit.ps1: could not find the 'a easily way' run to split"
}
$xml = $xml.Replace($old1, $new1)

# --- Change 2 --------------------------------------------------------
# After "The sub goal ... each color." paragraph: drop the trailing
# bookmark from that paragraph, then insert a blank paragraph, an
# "Identify Potential Solutions" heading, and the new solution
# paragraph; the bookmark is relocated into what used to be the final
# (empty) paragraph of the document.
$old2 = '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The sub goal is to get enough socks to make the selections of at least one matching pair and one matching pair of each color.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w14:paraId="75348CA2" w14:textId="77777777" w:rsidR="00A35249" w:rsidRPr="00A35249" w:rsidRDefault="00A35249"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'

$new2 = '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>The sub goal is to get enough socks to make the selections of at least one matching pair and one matching pair of each color.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Identify Potential Solutions</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>To me the possible solution would be to grab all 20 pairs while in the dark and then make the proper selections when you get to some light.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

if (-not $xml.Contains($old2)) {
    throw "edit.ps1: could not find the sock sub-goal / trailing-paragraph block"
}
$xml = $xml.Replace($old2, $new2)

# Write the whole story back.
$rng.InsertXML($xml)
